# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell updates: Row, Column, NewValue
$updates = @(
    @{ Row = 2; Col = 'D'; Value = '62.374.51' },
    @{ Row = 2; Col = 'E'; Value = '  -2.32%  ' },
    @{ Row = 3; Col = 'D'; Value = '2.623.07' },
    @{ Row = 3; Col = 'E'; Value = '  -4.07%  ' },
    @{ Row = 4; Col = 'E'; Value = '  -0.04%  ' },
    @{ Row = 5; Col = 'D'; Value = '551.73' },
    @{ Row = 5; Col = 'E'; Value = '  -2.27%  ' },
    @{ Row = 6; Col = 'D'; Value = '154.21' },
    @{ Row = 6; Col = 'E'; Value = '  -4.13%  ' },
    @{ Row = 7; Col = 'D'; Value = '0.999' },
    @{ Row = 7; Col = 'E'; Value = '  +0.04%  ' },
    @{ Row = 8; Col = 'E'; Value = '  -1.45%  ' },
    @{ Row = 9; Col = 'E'; Value = '  -4.00%  ' },
    @{ Row = 10; Col = 'E'; Value = '  -4.40%  ' },
    @{ Row = 11; Col = 'E'; Value = '  -1.43%  ' },
    @{ Row = 12; Col = 'D'; Value = '0.364' },
    @{ Row = 12; Col = 'E'; Value = '  -3.03%  ' },
    @{ Row = 13; Col = 'D'; Value = '3.085.11' },
    @{ Row = 13; Col = 'E'; Value = '  -4.05%  ' },
    @{ Row = 14; Col = 'E'; Value = '  -4.35%  ' },
    @{ Row = 15; Col = 'D'; Value = '62.224.80' },
    @{ Row = 15; Col = 'E'; Value = '  -2.30%  ' },
    @{ Row = 16; Col = 'E'; Value = '  -3.48%  ' },
    @{ Row = 17; Col = 'D'; Value = '2.625.10' },
    @{ Row = 17; Col = 'E'; Value = '  -3.91%  ' },
    @{ Row = 18; Col = 'D'; Value = '11.61' },
    @{ Row = 18; Col = 'E'; Value = '  -5.81%  ' },
    @{ Row = 19; Col = 'D'; Value = '4.52' },
    @{ Row = 19; Col = 'E'; Value = '  -4.03%  ' },
    @{ Row = 20; Col = 'D'; Value = '339.72' },
    @{ Row = 20; Col = 'E'; Value = '  -4.52%  ' },
    @{ Row = 21; Col = 'E'; Value = '  -8.01%  ' },
    @{ Row = 22; Col = 'D'; Value = '0.996' },
    @{ Row = 22; Col = 'E'; Value = '  -0.30%  ' },
    @{ Row = 23; Col = 'E'; Value = '  -4.06%  ' },
    @{ Row = 24; Col = 'D'; Value = '62.81' },
    @{ Row = 24; Col = 'E'; Value = '  -2.06%  ' },
    @{ Row = 25; Col = 'D'; Value = '0.167' },
    @{ Row = 25; Col = 'E'; Value = '  -1.24%  ' },
    @{ Row = 26; Col = 'E'; Value = '  -0.11%  ' },
    @{ Row = 27; Col = 'D'; Value = '8.02' },
    @{ Row = 27; Col = 'E'; Value = '  -3.67%  ' },
    @{ Row = 28; Col = 'D'; Value = '0.0₃0827' },
    @{ Row = 28; Col = 'E'; Value = '  -8.83%  ' },
    @{ Row = 29; Col = 'B'; Value = 'Fetch.AI' },
    @{ Row = 29; Col = 'C'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' },
    @{ Row = 29; Col = 'D'; Value = '1.33' },
    @{ Row = 29; Col = 'E'; Value = '  -3.15%  ' },
    @{ Row = 30; Col = 'B'; Value = 'Aptos' },
    @{ Row = 30; Col = 'C'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Row = 30; Col = 'D'; Value = '7.09' },
    @{ Row = 30; Col = 'E'; Value = '  -1.22%  ' },
    @{ Row = 31; Col = 'E'; Value = '  -4.13%  ' },
    @{ Row = 32; Col = 'B'; Value = 'USDe' },
    @{ Row = 32; Col = 'C'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' },
    @{ Row = 32; Col = 'D'; Value = '0.999' },
    @{ Row = 32; Col = 'E'; Value = '  +0.03%  ' },
    @{ Row = 33; Col = 'B'; Value = 'Monero' },
    @{ Row = 33; Col = 'C'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Row = 33; Col = 'D'; Value = '159.71' },
    @{ Row = 33; Col = 'E'; Value = '  -4.38%  ' },
    @{ Row = 34; Col = 'D'; Value = '4.71' },
    @{ Row = 34; Col = 'E'; Value = '  -4.39%  ' },
    @{ Row = 35; Col = 'E'; Value = '  -5.03%  ' },
    @{ Row = 36; Col = 'D'; Value = '19.21' },
    @{ Row = 36; Col = 'E'; Value = '  -4.14%  ' },
    @{ Row = 37; Col = 'E'; Value = '  -4.28%  ' },
    @{ Row = 38; Col = 'D'; Value = '335.39' },
    @{ Row = 38; Col = 'E'; Value = '  -2.58%  ' },
    @{ Row = 39; Col = 'D'; Value = '6.13' },
    @{ Row = 39; Col = 'E'; Value = '  -2.73%  ' },
    @{ Row = 40; Col = 'D'; Value = '0.892' },
    @{ Row = 40; Col = 'E'; Value = '  -8.19%  ' },
    @{ Row = 41; Col = 'D'; Value = '3.90' },
    @{ Row = 41; Col = 'E'; Value = '  -3.94%  ' },
    @{ Row = 42; Col = 'D'; Value = '37.69' },
    @{ Row = 42; Col = 'E'; Value = '  -2.47%  ' },
    @{ Row = 43; Col = 'D'; Value = '0.998' },
    @{ Row = 43; Col = 'E'; Value = '  -0.04%  ' },
    @{ Row = 44; Col = 'D'; Value = '20.33' },
    @{ Row = 44; Col = 'E'; Value = '  -6.50%  ' },
    @{ Row = 45; Col = 'D'; Value = '0.608' },
    @{ Row = 45; Col = 'E'; Value = '  -3.61%  ' },
    @{ Row = 46; Col = 'E'; Value = '  -0.68%  ' },
    @{ Row = 47; Col = 'D'; Value = '19.68' },
    @{ Row = 47; Col = 'E'; Value = '  -5.88%  ' },
    @{ Row = 48; Col = 'D'; Value = '0.0546' },
    @{ Row = 48; Col = 'E'; Value = '  -6.38%  ' },
    @{ Row = 49; Col = 'D'; Value = '0.0960' },
    @{ Row = 49; Col = 'E'; Value = '  -3.37%  ' },
    @{ Row = 50; Col = 'B'; Value = 'Maker' },
    @{ Row = 50; Col = 'C'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Row = 50; Col = 'D'; Value = '2.082.51' },
    @{ Row = 50; Col = 'E'; Value = '  -1.73%  ' },
    @{ Row = 51; Col = 'B'; Value = 'Aave' },
    @{ Row = 51; Col = 'C'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Row = 51; Col = 'D'; Value = '127.37' },
    @{ Row = 51; Col = 'E'; Value = '  -3.76%  ' }
)

# Regex that matches what Excel will auto-parse as a plain decimal number.
$numericPattern = '^-?\d+(\.\d+)?$'

foreach ($u in $updates) {
    $addr = "$($u.Col)$($u.Row)"
    $range = $ws.Range($addr)
    if ($u.Value -match $numericPattern) {
        # Force text storage so numeric-looking strings (e.g. "0.999") stay text,
        # matching the original inlineStr cell type, then restore default styling.
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}

Write-Host "Applied $($updates.Count) cell updates"
